$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 170, shifting existing rows 170-216 down to 171-217.
$ws.Rows.Item(170).Insert()

# New row 170 shares the static/categorical fields with the rest of this
# "Feria Lagunitas de Puerto Montt" / Pina block (rows 169 & 171 before insert).
$ws.Cells.Item(170, 1).Value = 4
$ws.Cells.Item(170, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(170, 3).Value = "Los Lagos"
$ws.Cells.Item(170, 4).Value = 44642
$ws.Cells.Item(170, 4).NumberFormat = $ws.Cells.Item(171, 4).NumberFormat
$ws.Cells.Item(170, 5).Value = 10
$ws.Cells.Item(170, 6).Value = "Fruta"
$ws.Cells.Item(170, 7).Value = 100108
$ws.Cells.Item(170, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(170, 9).Value = 100108005
$ws.Cells.Item(170, 10).Value = "Piña"
$ws.Cells.Item(170, 11).Value = "Caramelo"
$ws.Cells.Item(170, 12).Value = "Tercera"
$ws.Cells.Item(170, 13).Value = 200
$ws.Cells.Item(170, 14).Value = 17000
$ws.Cells.Item(170, 15).Value = 18000
$ws.Cells.Item(170, 16).Value = 17500
$ws.Cells.Item(170, 17).Value = "`$/caja 16 unidades"
$ws.Cells.Item(170, 18).Value = "Ecuador"
$ws.Cells.Item(170, 19).Value = 1094
$ws.Cells.Item(170, 20).Value = 16
